$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "Phường Tăng Nhơn Phú"
    3  = "Phường Chợ Lớn"
    4  = "Phường Tân Định"
    5  = "Phường Linh Xuân"
    6  = "Phường Tân Uyên"
    7  = "Phường Cầu Ông Lãnh"
    8  = "Phường Cầu Kiệu"
    9  = "Phường Tam Thắng"
    10 = "Phường Bình Lợi Trung"
    11 = "Xã Kim Long"
    12 = "Xã Xuyên Mộc"
    13 = "Xã Ngãi Giao"
    14 = "Phường Việt Hưng"
    15 = "Phường Thượng Cát"
    16 = "Phường Lĩnh Nam"
    17 = "Phường Yên Nghĩa"
    18 = "Phường Bồ Đề"
    19 = "Phường Kiến Hưng"
    20 = "Phường Ba Đình"
    21 = "Xã Liên Minh"
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}
